# Typhon_Profits workbook: refresh market-price-derived columns (H:N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit tables.
# Columns hold plain numbers (no formulas):
#   H currentAveragePrice      I currentAveragePriceNQ
#   J currentAveragePriceHQ    K LevePriceNQ
#   L LevePriceHQ              M LeveProfitNQ   N LeveProfitHQ
# Values below mirror a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8: On the Drip
$ws.Range("H8").Value = 133.33333
$ws.Range("J8").Value = 200
$ws.Range("L8").Value = 600
$ws.Range("N8").Value = -878
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 1204.0526
$ws.Range("I40").Value = 770
$ws.Range("J40").Value = 1686.3334
$ws.Range("K40").Value = 770
$ws.Range("L40").Value = 1686.3334
$ws.Range("M40").Value = -595
$ws.Range("N40").Value = -2036.3334
# Row 42: Eye of the Beholder
$ws.Range("H42").Value = 92.57143000000001
$ws.Range("I42").Value = 60
$ws.Range("J42").Value = 98
$ws.Range("K42").Value = 180
$ws.Range("L42").Value = 294
$ws.Range("M42").Value = 50
$ws.Range("N42").Value = -754
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 1200
$ws.Range("I43").Value = 1200
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1200
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1131
$ws.Range("N43").ClearContents()
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 463.9524
$ws.Range("I92").Value = 337.53333
$ws.Range("K92").Value = 337.53333
$ws.Range("M92").Value = 910.46667
# Row 98: The Dotted Line
$ws.Range("H98").Value = 598.8889
$ws.Range("I98").Value = 655.7143
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 655.7143
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 842.2857
$ws.Range("N98").Value = -3396
# Row 103: Let Loose the Juice
$ws.Range("H103").Value = 186
$ws.Range("J103").Value = 134
$ws.Range("L103").Value = 402
$ws.Range("N103").Value = -1574
# Row 106: Making Your Mark
$ws.Range("H106").Value = 8774116
$ws.Range("I106").Value = 12347287
$ws.Range("K106").Value = 12347287
$ws.Range("M106").Value = -12346656
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 4630687
$ws.Range("J112").Value = 5556666.5
$ws.Range("L112").Value = 16669999.5
$ws.Range("N112").Value = -16672215.5
# Row 122: Wishful Inking
$ws.Range("H122").Value = 598.8889
$ws.Range("I122").Value = 655.7143
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 1967.1429
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 482.8571000000002
$ws.Range("N122").Value = -6100
# Row 125: Body over Mind
$ws.Range("H125").Value = 938.6667
$ws.Range("I125").Value = 390
$ws.Range("K125").Value = 3510
$ws.Range("M125").Value = -1050
# Row 131: Mindful Study
$ws.Range("H131").Value = 1499.4595
$ws.Range("J131").Value = 2127.7778
$ws.Range("L131").Value = 6383.3334
$ws.Range("N131").Value = -16463.3334
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1580.625
$ws.Range("I137").Value = 1572.6666
$ws.Range("J137").Value = 1700
$ws.Range("K137").Value = 4717.9998
$ws.Range("L137").Value = 5100
$ws.Range("M137").Value = -2167.9998
$ws.Range("N137").Value = -10200

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6076.815
$ws.Range("I32").Value = 5222.68
$ws.Range("K32").Value = 5222.68
$ws.Range("M32").Value = -4935.68
# Row 35: Need for Mead
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -594
$ws.Range("N35").ClearContents()
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1876.3
$ws.Range("I102").Value = 1876.3
$ws.Range("K102").Value = 1876.3
$ws.Range("M102").Value = -254.3

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 37: That's Some Fine Grinding
$ws.Range("H37").Value = 1026
$ws.Range("I37").Value = 1026
$ws.Range("K37").Value = 1026
$ws.Range("M37").Value = -889

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 1199.5
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2574
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3530.25
$ws.Range("I31").Value = 2689.375
$ws.Range("J31").Value = 4202.95
$ws.Range("K31").Value = 2689.375
$ws.Range("L31").Value = 4202.95
$ws.Range("M31").Value = -2394.375
$ws.Range("N31").Value = -4792.95
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3530.25
$ws.Range("I34").Value = 2689.375
$ws.Range("J34").Value = 4202.95
$ws.Range("K34").Value = 2689.375
$ws.Range("L34").Value = 4202.95
$ws.Range("M34").Value = -2487.375
$ws.Range("N34").Value = -4606.95
# Row 107: Built to Last
$ws.Range("H107").Value = 983.5484
$ws.Range("I107").Value = 435.44446
$ws.Range("K107").Value = 435.44446
$ws.Range("M107").Value = 1484.55554
# Row 113: Patient Patients
$ws.Range("H113").Value = 1199.5
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 4463.875
$ws.Range("I132").Value = 3282.7273
$ws.Range("K132").Value = 9848.1819
$ws.Range("M132").Value = -7318.1819

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 57: The Egg Files
$ws.Range("H57").Value = 10000
$ws.Range("J57").Value = 10000
$ws.Range("L57").Value = 30000
$ws.Range("N57").Value = -31118
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 379.95
$ws.Range("J113").Value = 377.4375
$ws.Range("L113").Value = 1132.3125
$ws.Range("N113").Value = -5472.3125
# Row 129: Comfort Food
$ws.Range("H129").Value = 255944.75
$ws.Range("J129").Value = 464805.9
$ws.Range("L129").Value = 1394417.7
$ws.Range("N129").Value = -1404417.7
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 727.54
$ws.Range("I131").Value = 570
$ws.Range("J131").Value = 734.1042
$ws.Range("K131").Value = 1710
$ws.Range("L131").Value = 2202.3126
$ws.Range("M131").Value = 3330
$ws.Range("N131").Value = -12282.3126
# Row 139: Najoothie
$ws.Range("H139").Value = 2486.2188
$ws.Range("I139").Value = 1607.7778
$ws.Range("K139").Value = 4823.3334
$ws.Range("M139").Value = 316.6665999999996

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 16669015
$ws.Range("I102").Value = 22729658
$ws.Range("J102").Value = 2246
$ws.Range("K102").Value = 22729658
$ws.Range("L102").Value = 2246
$ws.Range("M102").Value = -22728036
$ws.Range("N102").Value = -5490
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2708
$ws.Range("I113").Value = 2074.5454
$ws.Range("J113").Value = 4450
$ws.Range("K113").Value = 2074.5454
$ws.Range("L113").Value = 4450
$ws.Range("M113").Value = 95.45460000000003
$ws.Range("N113").Value = -8790

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 100.916664
$ws.Range("J55").Value = 107.2
$ws.Range("L55").Value = 107.2
$ws.Range("N55").Value = -453.2
# Row 112: A Slippery Slope
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 1937.5385
$ws.Range("I136").Value = 1937.5385
$ws.Range("K136").Value = 5812.6155
$ws.Range("M136").Value = -3262.6155

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 61: Bundle Up, It's Odd out There
$ws.Range("H61").Value = 8051
$ws.Range("I61").Value = 8051
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8051
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7759
$ws.Range("N61").ClearContents()
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1185.68
$ws.Range("I126").Value = 1222.3636
$ws.Range("K126").Value = 3667.0908
$ws.Range("M126").Value = -1197.0908
# Row 130: Skill Cap
$ws.Range("H130").Value = 33996.43
$ws.Range("J130").Value = 33996.43
$ws.Range("L130").Value = 33996.43
$ws.Range("N130").Value = -44036.43
